# Add a new worksheet "strategy_id-6006" that is a duplicate of the
# existing "strategy_id-6004" sheet (same header row + same single data
# row), placed immediately after it in the workbook.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("strategy_id-6004")

# Copy the source sheet to a position right after itself.
$sourceSheet.Copy([System.Reflection.Missing]::Value, $sourceSheet)

# Excel names the freshly-copied sheet "strategy_id-6004 (2)"; rename it.
$newSheet = $wb.Worksheets.Item("strategy_id-6004 (2)")
$newSheet.Name = "strategy_id-6006"
